$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove hyperlinks from E2:E5 and reset their style to Normal ---
$ws.Hyperlinks.Delete()
$ws.Range("E2:E5").Style = "Normal"
$wb.Styles("Hyperlink").Delete()

# --- Update existing rows with new data ---
# Row 2 (Kiran)
$ws.Range("A2").Value = "Kiran"
$ws.Range("B2").Value = "refinitiv"
$ws.Range("C2").Value = "ny4v3"
$ws.Range("D2").Value = "Email"
$ws.Range("E2").Value = "kiran@gmail.com"
$ws.Range("F2").Value = "RIC"
$ws.Range("G2").Value = $true

# Row 3 (Rahul)
$ws.Range("A3").Value = "Rahul"
$ws.Range("B3").Value = "refinitiv"
$ws.Range("D3").Value = "Jira"
$ws.Range("E3").Value = "rahul@gmail.com"
$ws.Range("F3").Value = "PDP"
$ws.Range("G3").Value = $false

# Row 4 (Rajendra)
$ws.Range("A4").Value = "Rajendra"
$ws.Range("B4").Value = "refinitiv"
$ws.Range("C4").Value = "ny4v3"
$ws.Range("D4").Value = "Jira"
$ws.Range("E4").Value = "rajendra@gmail.com"
$ws.Range("F4").Value = "PERM"
$ws.Range("G4").Value = $true

# Row 5 (Mahendra)
$ws.Range("A5").Value = "Mahendra"
$ws.Range("B5").Value = "bloomberg"
$ws.Range("D5").Value = "Email"
$ws.Range("E5").Value = "mahendra@gmail.com"
$ws.Range("G5").Value = $true

# --- New row 6 (Pradip) ---
$ws.Range("A6").Value = "Pradip"
$ws.Range("B6").Value = "prodcut"
$ws.Range("C6").Value = " north_america"
$ws.Range("D6").Value = "jira"
$ws.Range("E6").Value = "pradip@gmail.com"
$ws.Range("F6").Value = "PE"
$ws.Range("G6").Value = $false

# --- Column widths (closest achievable values given this engine's width quantization) ---
$ws.Columns.Item(2).ColumnWidth = 18.5
$ws.Columns.Item(3).ColumnWidth = 22.0
$ws.Columns.Item(4).ColumnWidth = 14.5
$ws.Columns.Item(5).ColumnWidth = 37.65
$ws.Columns.Item(6).ColumnWidth = 9.0

# --- Selection ---
$ws.Range("F6").Select() | Out-Null
